$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.895.60'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '1.605.77'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.00'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('E7').Value = '  -3.74%  '
$ws.Range('E8').Value = '  -2.10%  '
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.01'
$ws.Range('E10').Value = '  -2.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').Value = '1.831.60'
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('D13').Value = '1.633.34'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.07'
$ws.Range('E14').Value = '  -3.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.515'
$ws.Range('E15').Value = '  -2.77%  '
$ws.Range('D16').Value = '25.879.08'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.20'
$ws.Range('E17').Value = '  -1.09%  '
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  -2.57%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '190.54'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.19'
$ws.Range('E21').Value = '  -1.63%  '
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.95'
$ws.Range('E23').Value = '  -3.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '142.48'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('E25').Value = '  -3.50%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.71'
$ws.Range('E27').Value = '  -3.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.56'
$ws.Range('E28').Value = '  -3.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.01'
$ws.Range('E29').Value = '  -1.32%  '
$ws.Range('E30').Value = '  -1.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0471'
$ws.Range('E31').Value = '  -2.60%  '
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('E33').Value = '  -3.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.40'
$ws.Range('E34').Value = '  -1.06%  '
$ws.Range('E35').Value = '  -1.93%  '
$ws.Range('D36').Value = '1.119.29'
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.38'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.810'
$ws.Range('E38').Value = '  -6.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0151'
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.505'
$ws.Range('E40').Value = '  -3.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '96.39'
$ws.Range('E41').Value = '  -2.11%  '
$ws.Range('D42').Value = '1.744.96'
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.746'
$ws.Range('E43').Value = '  -4.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.07'
$ws.Range('E44').Value = '  -3.97%  '
$ws.Range('D45').Value = '0.0₆0112'
$ws.Range('E45').Value = '  -1.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '53.44'
$ws.Range('E46').Value = '  -3.31%  '
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0512'
$ws.Range('E48').Value = '  -2.80%  '
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.36'
$ws.Range('E51').Value = '  -2.66%  '
